$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the three new "Akurasi" columns (one per Indikator/Data block).
#    Columns("C") insert -> new blank col C, old C..F shift to D..G
#    Columns("F") insert -> lands right after the (now shifted) Data2 column E
#    Columns("I") insert -> lands right after the (now shifted) Data3 column H
# ---------------------------------------------------------------------------
$ws.Columns("C").Insert()
$ws.Columns("F").Insert()
$ws.Columns("I").Insert()

# ---------------------------------------------------------------------------
# 2) Header row 2: label the new columns "Akurasi", matching style of Data col.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Akurasi"
$ws.Range("F2").Value = "Akurasi"
$ws.Range("I2").Value = "Akurasi"

# ---------------------------------------------------------------------------
# 3) Number format + centered alignment for the new Akurasi columns. Each
#    block only has data through its own last populated row, so format only
#    those cells (no styling on the otherwise-untouched rows below).
# ---------------------------------------------------------------------------
$akurasiRanges = @("C3:C6", "C10", "F3:F5", "F10", "I3:I10")
foreach ($r in $akurasiRanges) {
    $ws.Range($r).NumberFormat = "0.00%"
    $ws.Range($r).HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 4) Data values - block 1 (Directory Information)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 252
$ws.Range("C6").Value = 0.5187007874015748

# ---------------------------------------------------------------------------
# 5) Data values - block 2 (Educational Information)
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# ---------------------------------------------------------------------------
# 6) Data values - block 3 (Personally Identifiable)
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = 252
$ws.Range("I3").Value = 0.5187007874015748
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("H9").Value = 214
$ws.Range("I9").Value = 0.421259842519685

# ---------------------------------------------------------------------------
# 7) Totals row (row 10): add AVERAGE formulas for the new Akurasi columns,
#    fix the block-3 Data total (SUM) range to extend through row 9.
# ---------------------------------------------------------------------------
$ws.Range("C10").Formula = "=AVERAGE(C1:C7)"
$ws.Range("F10").Formula = "=AVERAGE(F1:F7)"
$ws.Range("H10").Formula = "=SUM(H1:H9)"
$ws.Range("I10").Formula = "=AVERAGE(I1:I9)"

# ---------------------------------------------------------------------------
# 8) Fix up the conditional-formatting ranges before the row delete below -
#    set the new formula text first, THEN move the applies-to range, so the
#    persisted rule keeps its original rule-type/dxf but points at the new
#    range and formula.
# ---------------------------------------------------------------------------
$fcTotal = $ws.Range("A1:F10").FormatConditions
for ($i = 1; $i -le $fcTotal.Count; $i++) {
    $fcTotal.Item($i).ModifyAppliesToRange($ws.Range("A1:I10"))
}

$fcFooter = $ws.Range("A12:F12").FormatConditions.Item(1)
$fcFooter.Formula1 = "=LEN(TRIM(A11))>0"
$fcFooter.ModifyAppliesToRange($ws.Range("A11:I11"))

# ---------------------------------------------------------------------------
# 9) Update the footer label text, then delete the now-empty row 11 so the
#    footer row (old row 12) shifts up to row 11.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Akurasi Pengujian = 8.80%"
$ws.Rows(11).Delete()

# ---------------------------------------------------------------------------
# 10) Re-merge the header band cells (inserting columns inside the old merged
#     range does not auto-widen it) and the footer row merge.
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").MergeCells = $true
$ws.Range("D1:F1").MergeCells = $true
$ws.Range("G1:I1").MergeCells = $true
$ws.Range("A11:I11").MergeCells = $true
